# Fix the weird ENVO term in the AgDaFAIR potato field imaging template.
#
# The "Tags Term Accession Number" for the "agricultural field" tag used a
# malformed OBO PURL; swap it for the proper bioregistry.io link (and make
# it a real clickable hyperlink, like the other bioregistry.io links in the
# row above). Also fix the matching "Tags Term Source REF" abbreviation from
# the lowercase "envo" to the canonical "ENVO".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# D14: Tags Term Accession Number (agricultural field / ENVO)
$ws.Range("D14").Value = "https://bioregistry.io/ENVO:00000114"
$d14 = $ws.Cells.Item(14, 4)
$ws.Hyperlinks.Add($d14, "https://bioregistry.io/envo:00000114", "", "https://bioregistry.io/envo:00000114")

# D15: Tags Term Source REF (agricultural field / ENVO)
$ws.Range("D15").Value = "ENVO"

# Cosmetic follow-up: widen the columns to fit the (now longer) content and
# move the selection, matching the author's final view of the sheet.
$ws.Columns.Item(1).ColumnWidth = 34.290625
$ws.Columns.Item(2).ColumnWidth = 83.57578125
$ws.Columns.Item(3).ColumnWidth = 135.57578125
$ws.Columns.Item(4).ColumnWidth = 113.15

$ws.Range("C22").Select()
